$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (YHL), shifting NIEL (and below) down one row
$ws.Rows.Item(3).Insert()

# Row 2: KHRYZ - update date
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-02-22"
$ws.Range("B2").Style = "Normal"

# Row 3: new entry YHL
$ws.Range("A3").Value = "YHL"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2025-02-22"
$ws.Range("B3").Style = "Normal"

# Row 4: NIEL (shifted down from row 3) - update date
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2025-02-22"
$ws.Range("B4").Style = "Normal"

# Row 5: new entry KIM
$ws.Range("A5").Value = "KIM"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2025-02-22"
$ws.Range("B5").Style = "Normal"

# Row 6: new entry 2029123_KIMI ROKKU
$ws.Range("A6").Value = "2029123_KIMI ROKKU"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2025-02-22"
$ws.Range("B6").Style = "Normal"
